$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column (D) keeps its text formatting so values
# like "521.96" are not auto-converted to numbers by Excel.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = "59.363.81"
$ws.Range("E2").Value = "  +4.15%  "
$ws.Range("D3").Value = "2.599.00"
$ws.Range("E3").Value = "  +2.44%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "521.96"
$ws.Range("E5").Value = "  +1.98%  "
$ws.Range("D6").Value = "140.85"
$ws.Range("E6").Value = "  +1.23%  "
$ws.Range("E7").Value = "  -0.37%  "
$ws.Range("E8").Value = "  +1.85%  "
$ws.Range("D9").Value = "2.623.10"
$ws.Range("E9").Value = "  +3.17%  "
$ws.Range("D10").Value = "6.53"
$ws.Range("E10").Value = "  +1.05%  "
$ws.Range("E11").Value = "  +2.20%  "
$ws.Range("E12").Value = "  +2.57%  "
$ws.Range("E13").Value = "  +1.99%  "
$ws.Range("D14").Value = "3.060.97"
$ws.Range("E14").Value = "  +2.59%  "
$ws.Range("D15").Value = "59.369.23"
$ws.Range("E15").Value = "  +4.11%  "
$ws.Range("E16").Value = "  +2.17%  "
$ws.Range("D17").Value = "2.617.50"
$ws.Range("E17").Value = "  +3.41%  "
$ws.Range("E18").Value = "  +0.55%  "
$ws.Range("D19").Value = "338.89"
$ws.Range("E19").Value = "  +2.37%  "
$ws.Range("D20").Value = "4.33"
$ws.Range("E20").Value = "  +1.38%  "
$ws.Range("E21").Value = "  +1.35%  "
$ws.Range("D22").Value = "6.51"
$ws.Range("E22").Value = "  +6.60%  "
$ws.Range("E23").Value = "  -0.25%  "
$ws.Range("D24").Value = "66.26"
$ws.Range("E24").Value = "  +3.43%  "
$ws.Range("E25").Value = "  +1.57%  "
$ws.Range("E26").Value = "  +1.24%  "
$ws.Range("E27").Value = "  -0.69%  "
$ws.Range("E28").Value = "  +2.18%  "
$ws.Range("E30").Value = "  -2.66%  "
$ws.Range("D31").Value = "5.94"
$ws.Range("E31").Value = "  -5.16%  "
$ws.Range("D32").Value = "18.79"
$ws.Range("E32").Value = "  +2.22%  "
$ws.Range("E33").Value = "  +1.19%  "
$ws.Range("D34").Value = "149.02"
$ws.Range("E34").Value = "  +0.27%  "
$ws.Range("E35").Value = "  +1.40%  "
$ws.Range("E36").Value = "  +0.38%  "
$ws.Range("E37").Value = "  +1.97%  "
$ws.Range("E38").Value = "  +4.51%  "
$ws.Range("D39").Value = "0.834"
$ws.Range("E39").Value = "  +2.01%  "
$ws.Range("D40").Value = "0.821"
$ws.Range("E40").Value = "  -2.09%  "
$ws.Range("E41").Value = "  +2.36%  "
$ws.Range("E42").Value = "  -0.59%  "
$ws.Range("D43").Value = "275.83"
$ws.Range("E43").Value = "  +7.47%  "
$ws.Range("D44").Value = "10.72"
$ws.Range("E44").Value = "  +1.18%  "
$ws.Range("D45").Value = "0.592"
$ws.Range("E45").Value = "  +3.05%  "
$ws.Range("E46").Value = "  +0.33%  "
$ws.Range("E47").Value = "  +0.24%  "
$ws.Range("D48").Value = "18.64"
$ws.Range("E48").Value = "  +1.19%  "
$ws.Range("D49").Value = "1.985.04"
$ws.Range("E49").Value = "  +0.95%  "
$ws.Range("E50").Value = "  +2.87%  "
$ws.Range("E51").Value = "  +0.13%  "

# Restore the original (default) cell style now that the text values
# are safely stored, so no extraneous formatting is introduced.
$priceRange.Style = "Normal"
